$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.318692803382874
$ws.Range("B1").Value = 3.286046028137207
$ws.Range("C1").Value = 5.640407562255859
$ws.Range("D1").Value = 1.717574119567871
$ws.Range("E1").Value = 1.005552411079407
